$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (index 1): bump "想去人数" (column F) counts on a number of rows.
# ---------------------------------------------------------------------------
$wsExhibit = $wb.Worksheets.Item(1)
$exhibitUpdates = @{
    3  = 401
    5  = 39
    6  = 1257
    7  = 463
    9  = 215
    11 = 188
    12 = 1064
    13 = 7
    14 = 276
    15 = 208
    16 = 1544
    17 = 566
    18 = 239
    19 = 361
    21 = 864
    22 = 1171
    25 = 2693
    26 = 1480
    28 = 57
    29 = 468
    30 = 820
    31 = 1354
    33 = 1438
    36 = 797
    37 = 658
    38 = 696
    39 = 887
    41 = 265
}
foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Cells.Item($row, 6).Value = $exhibitUpdates[$row]
}

# ---------------------------------------------------------------------------
# Sheet "演出" (index 2):
#   * bump F15
#   * insert a new event row at row 22 (a Beyond tribute concert), which
#     pushes the three rows that used to be 22-24 down to 23-25
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item(2)

$wsShow.Cells.Item(15, 6).Value = 665

# Insert the new row above the current row 22.
$wsShow.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new concert's data.
$a22 = $wsShow.Range("A22")
$a22.Value = 21
$a22.Font.Bold = $true
$a22.HorizontalAlignment = -4108
$a22.VerticalAlignment = -4160
$a22.Borders.LineStyle = 1

$b22 = $wsShow.Range("B22")
$b22.NumberFormat = "@"
$b22.Value = "2024-08-10"
$wsShow.Range("C22").Value = "杭州·【七夕巨献·早鸟6折】真的爱你”致敬Beyond·黄家驹31周年演唱会·630乐团再现91殿堂级演出"
$wsShow.Range("D22").Value = "湖墅南路136-138号 浙话艺术剧院"
$wsShow.Range("E22").Value = "2024.08.10 19:30-08.10 21:30"
$wsShow.Cells.Item(22, 6).Value = 0
$wsShow.Cells.Item(22, 7).Value = 60
$wsShow.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=85333"
$wsShow.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202405/uYt32zt21715221330023.jpeg"

# Fix up the numbering (column A) and the interest count (column F) on the
# rows that were shifted down by the insert.
$wsShow.Cells.Item(23, 1).Value = 22
$wsShow.Cells.Item(23, 6).Value = 23

$wsShow.Cells.Item(24, 1).Value = 23

$wsShow.Cells.Item(25, 1).Value = 24

# ---------------------------------------------------------------------------
# Sheet "全部类型" (index 4): same kind of "想去人数" bumps as 展览, mirrored
# at different row offsets (this sheet is not touched by the 演出 row
# insertion above).
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item(4)
$allUpdates = @{
    4  = 401
    6  = 39
    9  = 1257
    10 = 463
    12 = 215
    14 = 188
    15 = 1064
    16 = 276
    18 = 208
    19 = 1544
    20 = 566
    21 = 239
    22 = 361
    25 = 1171
    26 = 2693
    28 = 1480
    31 = 57
    34 = 468
    35 = 820
    36 = 1354
    40 = 1438
    41 = 797
    42 = 658
    43 = 696
    44 = 887
    46 = 23
    48 = 265
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
